$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column D: header + values
$ws.Range("D1").Value = "target01"
$ws.Range("D2").Value = 1
$ws.Range("D3").Value = 1
$ws.Range("D4").Value = 1
$ws.Range("D5").Value = 0
$ws.Range("D6").Value = 0
$ws.Range("D7").Value = 0
$ws.Range("D8").Value = 0

# Update row 8 values in A and B
$ws.Range("A8").Value = -3
$ws.Range("B8").Value = 3

# Update selection to match the target view
$ws.Range("F4").Select()
